# Applies the "Checks for empty data" commit to geolocation_settings.xlsx:
#  - inserts a "log" column (after lightFile) and a "doTwilights" column
#    (after maxLightInt)
#  - removes the "deplAsk" column
#  - inserts "kernelStart" / "kernelEnd" columns (after createKernel)
#  - flips include -> TRUE and keepCalibPoints -> FALSE on the data row
#  - re-formats several date columns on the data row to the
#    mm-dd-yy / m/d/yy h:mm builtin formats
#  - adds three blank, pre-formatted template rows below the data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- 1. column structure -------------------------------------------------

$ws.Columns("D").Insert()
$ws.Columns("G").Insert()
$ws.Columns("AD").Delete()
$ws.Columns("AX").Insert()
$ws.Columns("AY").Insert()

$ws.Range("D1").Value = "log"
$ws.Range("G1").Value = "doTwilights"
$ws.Range("AX1").Value = "kernelStart"
$ws.Range("AY1").Value = "kernelEnd"

# ---- 2. data-row (row 2) value changes -----------------------------------

$ws.Range("D2").Value = $false
$ws.Range("G2").Value = $false

$ws.Range("B2").Value = $true
$ws.Range("AB2").Value = $false

# ---- 3. data-row (row 2) number-format changes ---------------------------

$ws.Range("I2:J2").NumberFormat = "mm-dd-yy"
$ws.Range("L2:M2").NumberFormat = "mm-dd-yy"
$ws.Range("O2:R2").NumberFormat = "mm-dd-yy"
$ws.Range("U2:V2").NumberFormat = "m/d/yy h:mm"

$ws.Range("AX2:AY2").NumberFormat = "yyyy\-mm\-dd;@"

# ---- 4. new blank template rows (3-5) ------------------------------------

$ws.Range("I3:J3").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("L3:M3").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("O3:R3").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("U3:V3").NumberFormat = "yyyy\-mm\-dd\ h:mm"
$ws.Range("AX3:AY3").NumberFormat = "yyyy\-mm\-dd;@"

$ws.Range("I4:J4").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("L4:M4").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("O4:R4").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("U4:V4").NumberFormat = "yyyy\-mm\-dd\ h:mm"
$ws.Range("AX4:AY4").NumberFormat = "yyyy\-mm\-dd;@"

$ws.Range("I5:J5").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("L5:M5").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("O5:R5").NumberFormat = "yyyy\-mm\-dd;@"
$ws.Range("U5:V5").NumberFormat = "yyyy\-mm\-dd\ h:mm"
$ws.Range("AX5:AY5").NumberFormat = "yyyy\-mm\-dd;@"

# ---- 5. sheet view tweaks --------------------------------------------------

$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 45
$ws.Range("BG6").Select()
